$d = $word.ActiveDocument

# --- Change 1 ---------------------------------------------------------
# The paragraph holding the inline picture with wp14:anchorId="10B2248F"
# is missing <w:rPr><w:noProof/></w:rPr> on its run. Find that paragraph
# (it's the one whose run contains an InlineShape but currently has no
# rPr at all) and turn on NoProofing, which Word represents as
# <w:noProof/> inside the run's rPr.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.InlineShapes.Count -gt 0 -and $p.Range.NoProofing -eq $false) {
        $target = $p
        break
    }
}
if ($target -ne $null) {
    $target.Range.NoProofing = $true
}

# --- Change 2 ---------------------------------------------------------
# Fix the typo "distncia" -> "distancia" in a way that leaves the
# sentence split across three runs, exactly as a manual retype of the
# missing "a" would: "...la dist" | "a" | "ncia entre tramos...".
$find = $d.Content
$find.Find.ClearFormatting()
$found = $find.Find.Execute("distncia", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found -eq $true) {
    $full = $find.Duplicate
    $splitPos = $full.Start + 4   # after "dist"

    $para = $full.Paragraphs(1)
    $paraRange = $para.Range
    $paraEnd = $paraRange.End - 1   # exclude the trailing paragraph mark

    $before = $d.Range($paraRange.Start, $splitPos).Text
    $after = $d.Range($splitPos, $paraEnd).Text

    $beforeEsc = $before.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    $afterEsc = $after.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>' + $beforeEsc + '</w:t></w:r><w:r><w:t>a</w:t></w:r><w:r><w:t>' + $afterEsc + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $target2 = $d.Range($paraRange.Start, $paraEnd)
    $target2.InsertXML($xml)
}

Write-Output "done"
